$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '26.900.72'
$ws.Cells.Item(2, 5).Value = '  +1.78%  '
$ws.Cells.Item(3, 4).Value = '1.727.23'
$ws.Cells.Item(3, 5).Value = '  +0.15%  '
$ws.Cells.Item(4, 4).Value = '0.9971'
$ws.Cells.Item(4, 5).Value = '  -0.22%  '
$ws.Cells.Item(5, 4).Value = '242.00'
$ws.Cells.Item(5, 5).Value = '  -0.47%  '
$ws.Cells.Item(6, 4).Value = '0.9974'
$ws.Cells.Item(6, 5).Value = '  -0.22%  '
$ws.Cells.Item(7, 4).Value = '0.4891'
$ws.Cells.Item(7, 5).Value = '  -0.52%  '
$ws.Cells.Item(8, 4).Value = '0.2593'
$ws.Cells.Item(8, 5).Value = '  -0.87%  '
$ws.Cells.Item(9, 4).Value = '0.06212'
$ws.Cells.Item(9, 5).Value = '  +0.07%  '
$ws.Cells.Item(10, 4).Value = '1.728.36'
$ws.Cells.Item(10, 5).Value = '  +0.31%  '
$ws.Cells.Item(11, 4).Value = '16.04'
$ws.Cells.Item(11, 5).Value = '  +3.70%  '
$ws.Cells.Item(12, 4).Value = '0.06894'
$ws.Cells.Item(12, 5).Value = '  -1.55%  '
$ws.Cells.Item(13, 4).Value = '0.6086'
$ws.Cells.Item(13, 5).Value = '  +1.46%  '
$ws.Cells.Item(14, 4).Value = '4.483'
$ws.Cells.Item(14, 5).Value = '  -1.85%  '
$ws.Cells.Item(15, 4).Value = '77.30'
$ws.Cells.Item(15, 5).Value = '  +0.04%  '
$ws.Cells.Item(16, 4).Value = '0.9975'
$ws.Cells.Item(16, 5).Value = '  -0.22%  '
$ws.Cells.Item(17, 4).Value = '26.656.58'
$ws.Cells.Item(18, 5).Value = '  -0.30%  '
$ws.Cells.Item(19, 4).Value = '0.000007185'
$ws.Cells.Item(19, 5).Value = '  +0.27%  '
$ws.Cells.Item(20, 5).Value = '  +0.83%  '
$ws.Cells.Item(21, 4).Value = '1.952.58'
$ws.Cells.Item(21, 5).Value = '  +0.59%  '
$ws.Cells.Item(22, 4).Value = '4.419'
$ws.Cells.Item(22, 5).Value = '  -1.52%  '
$ws.Cells.Item(23, 4).Value = '8.583'
$ws.Cells.Item(23, 5).Value = '  -0.03%  '
$ws.Cells.Item(24, 4).Value = '5.092'
$ws.Cells.Item(24, 5).Value = '  -1.42%  '
$ws.Cells.Item(25, 4).Value = '138.62'
$ws.Cells.Item(25, 5).Value = '  +0.33%  '
$ws.Cells.Item(26, 4).Value = '15.31'
$ws.Cells.Item(26, 5).Value = '  +0.31%  '
$ws.Cells.Item(27, 5).Value = '  +3.48%  '
$ws.Cells.Item(28, 4).Value = '106.30'
$ws.Cells.Item(28, 5).Value = '  -0.66%  '
$ws.Cells.Item(29, 4).Value = '1.379'
$ws.Cells.Item(29, 5).Value = '  -1.29%  '
$ws.Cells.Item(30, 4).Value = '3.953'
$ws.Cells.Item(30, 5).Value = '  +0.11%  '
$ws.Cells.Item(31, 4).Value = '0.07999'
$ws.Cells.Item(31, 5).Value = '  +0.33%  '
$ws.Cells.Item(32, 4).Value = '3.688'
$ws.Cells.Item(32, 5).Value = '  +0.30%  '
$ws.Cells.Item(33, 4).Value = '0.04523'
$ws.Cells.Item(33, 5).Value = '  -0.38%  '
$ws.Cells.Item(34, 4).Value = '0.9961'
$ws.Cells.Item(34, 5).Value = '  -0.28%  '
$ws.Cells.Item(35, 5).Value = '  -0.18%  '
$ws.Cells.Item(36, 4).Value = '1.012'
$ws.Cells.Item(36, 5).Value = '  +1.60%  '
$ws.Cells.Item(37, 4).Value = '0.6249'
$ws.Cells.Item(37, 5).Value = '  -0.22%  '
$ws.Cells.Item(38, 4).Value = '0.9375'
$ws.Cells.Item(38, 5).Value = '  +1.11%  '
$ws.Cells.Item(39, 4).Value = '2.055'
$ws.Cells.Item(39, 5).Value = '  +4.72%  '
$ws.Cells.Item(40, 4).Value = '2.460'
$ws.Cells.Item(40, 5).Value = '  +2.96%  '
$ws.Cells.Item(41, 4).Value = '0.9966'
$ws.Cells.Item(41, 5).Value = '  -0.30%  '
$ws.Cells.Item(42, 2).Value = 'VeChain'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(42, 4).Value = '0.01499'
$ws.Cells.Item(42, 5).Value = '  +0.87%  '
$ws.Cells.Item(43, 2).Value = 'FraxShare'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(43, 4).Value = '5.654'
$ws.Cells.Item(43, 5).Value = '  +5.99%  '
$ws.Cells.Item(44, 4).Value = '99.63'
$ws.Cells.Item(44, 5).Value = '  -0.20%  '
$ws.Cells.Item(45, 4).Value = '0.3852'
$ws.Cells.Item(45, 5).Value = '  +0.06%  '
$ws.Cells.Item(46, 4).Value = '6.851'
$ws.Cells.Item(46, 5).Value = '  +1.24%  '
$ws.Cells.Item(47, 4).Value = '0.1162'
$ws.Cells.Item(47, 5).Value = '  -0.61%  '
$ws.Cells.Item(48, 5).Value = '  +0.65%  '
$ws.Cells.Item(49, 4).Value = '7.898'
$ws.Cells.Item(49, 5).Value = '  +2.11%  '
$ws.Cells.Item(50, 4).Value = '30.20'
$ws.Cells.Item(50, 5).Value = '  +0.18%  '
$ws.Cells.Item(51, 4).Value = '51.65'
$ws.Cells.Item(51, 5).Value = '  +1.42%  '
